# Task #1080: rename TEK_ID column header to TEK, and rename the
# tek_period_start_year / tek_period_end_year headers to
# period_start_year / period_end_year. Also move the active selection
# from D12 to A2, matching the author's final saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) labels.
$ws.Range("A1").Value = "TEK"
$ws.Range("C1").Value = "period_start_year"
$ws.Range("D1").Value = "period_end_year"

# Move the selection to A2, like in the saved workbook.
$ws.Range("A2").Select() | Out-Null
